$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values in column D (speed/buff-time column)
$ws.Range("D2").Value = 0.1
$ws.Range("D3").Value = 1
$ws.Range("D5").Value = 1

# Update the active selection on the sheet to G7
$ws.Activate()
$ws.Range("G7").Select()
